# Updated automation test cases
# Adds a "${Status}" column (C) to the login data sheet, marking each
# set of credentials as "Invalid" or "Valid".
#
# NOTE: literal "${...}" text must be written with single-quoted strings
# so PowerShell does not try to interpolate ${Status} as a variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1, styled the same way as the existing A1/B1 headers
# (yellow fill).
$ws.Range("C1").Value = '${Status}'
$ws.Range("C1").Interior.Color = $ws.Range("A1").Interior.Color

# New status values for each existing data row.
$ws.Range("C2").Value = 'Invalid'
$ws.Range("C3").Value = 'Invalid'
$ws.Range("C4").Value = 'Valid'

# Move the active selection, matching the saved view state.
[void]$ws.Range("A6").Select()
